# Updated cryptos list with GitHub Actions
# Refreshes per-row Price (column D) and Volume(1h) (column E) text values,
# and swaps the Monero / PEPE rows (49-50), matching the source data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '63.638.09'
$ws.Range("E2").Value = '  +3.95%  '
# Row 3
$ws.Range("D3").Value = '3.072.12'
$ws.Range("E3").Value = '  +2.66%  '
# Row 4
$ws.Range("E4").Value = '  -0.03%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '552.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.79%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.08%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '
# Row 8
$ws.Range("D8").Value = '3.068.95'
$ws.Range("E8").Value = '  +2.62%  '
# Row 9
$ws.Range("E9").Value = '  +1.40%  '
# Row 10
$ws.Range("E10").Value = '  +6.64%  '
# Row 11
$ws.Range("E11").Value = '  +3.32%  '
# Row 12
$ws.Range("E12").Value = '  +1.87%  '
# Row 13
$ws.Range("E13").Value = '  +3.29%  '
# Row 14
$ws.Range("E14").Value = '  +2.56%  '
# Row 15
$ws.Range("D15").Value = '3.565.96'
$ws.Range("E15").Value = '  +2.38%  '
# Row 16
$ws.Range("D16").Value = '63.567.69'
$ws.Range("E16").Value = '  +3.80%  '
# Row 17
$ws.Range("D17").Value = '3.069.63'
$ws.Range("E17").Value = '  +2.56%  '
# Row 18
$ws.Range("E18").Value = '  -0.92%  '
# Row 19
$ws.Range("E19").Value = '  +2.26%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '486.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.24%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.97%  '
# Row 22
$ws.Range("E22").Value = '  +0.31%  '
# Row 23
$ws.Range("E23").Value = '  +4.90%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.58%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.30%  '
# Row 26
$ws.Range("E26").Value = '  +0.02%  '
# Row 27
$ws.Range("E27").Value = '  +3.12%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.31%  '
# Row 29
$ws.Range("E29").Value = '  +6.39%  '
# Row 30
$ws.Range("E30").Value = '  +0.08%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.28'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.62%  '
# Row 32
$ws.Range("E32").Value = '  +0.99%  '
# Row 33
$ws.Range("E33").Value = '  +8.35%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.69'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.92%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.50'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.16%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.04%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '465.91'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.07%  '
# Row 38
$ws.Range("E38").Value = '  +4.33%  '
# Row 39
$ws.Range("E39").Value = '  +3.07%  '
# Row 40
$ws.Range("D40").Value = '3.067.37'
$ws.Range("E40").Value = '  -3.14%  '
# Row 41
$ws.Range("E41").Value = '  +1.32%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.25'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.43%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.58'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.74%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '28.10'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.13%  '
# Row 45
$ws.Range("E45").Value = '  +4.55%  '
# Row 46
$ws.Range("E46").Value = '  -0.09%  '
# Row 47
$ws.Range("E47").Value = '  +2.80%  '
# Row 49
$ws.Range("B49").Value = 'PEPE'
$ws.Range("C49").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D49").Value = '0.0₃0514'
$ws.Range("E49").Value = '  +3.35%  '
# Row 50
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '117.29'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.53%  '
# Row 51
$ws.Range("E51").Value = '  +3.76%  '
